# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    cloning the layout/style of the "2021-Q4" sheet, and fill in the new
#    fund-holding figures for 2022-Q1.
# 2. Update the "总计" summary sheet: insert a new row for "2022-Q1" at the
#    top of the data (row 2) and shift the existing history down, updating
#    the running index column (A) and the 2021-Q4 market-value figure.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet, positioned right before "总计" ---
# Clone the most recent quarter sheet (keeps its layout/formatting/sheetPr
# intact) and drop the copy in right before "总计".
# NOTE: the worksheet reference passed as the "Before" sheet gets rebound to
# the newly-inserted sheet once Copy() runs, so re-fetch "总计" by name
# afterwards rather than reusing the handle used for placement.
$q4Sheet.Copy($wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Item("总计")

# Overwrite the cloned data row with the 2022-Q1 figures.
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "'010010"
$newSheet.Cells.Item(2, 3).Value = "国投瑞银港股通6个月定期开放股票"
$newSheet.Cells.Item(2, 4).Value = "'8.09"
$newSheet.Cells.Item(2, 5).Value = "'93.58"
$newSheet.Cells.Item(2, 6).Value = "'4.87"
$newSheet.Cells.Item(2, 7).Value = "'0.3940"
$newSheet.Cells.Item(2, 8).Value = 7

# --- 2. Update the "总计" sheet with the new 2022-Q1 row ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Give the new index cell (A2) the same style as the other index cells.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.39

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 1
$totalSheet.Cells.Item(3, 4).Value = 0.39

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(4, 3).Value = 1
$totalSheet.Cells.Item(4, 4).Value = 0.44

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(5, 3).Value = 1
$totalSheet.Cells.Item(5, 4).Value = 0.45

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(6, 3).Value = 3
$totalSheet.Cells.Item(6, 4).Value = 2.48

$totalSheet.Cells.Item(7, 1).Value = 5
$totalSheet.Cells.Item(7, 2).Value = "2020-Q4"
$totalSheet.Cells.Item(7, 3).Value = 1
$totalSheet.Cells.Item(7, 4).Value = 0.35

# Restore the originally active sheet/selection (workbook opened on the
# first sheet before this edit).
$firstSheet = $wb.Worksheets.Item(1)
$firstSheet.Activate()
$firstSheet.Range("A1").Select()
